$d = $word.ActiveDocument

# The document ends with a level-2 ("ilvl=2") bullet ("...similar properties of
# inputDict"). We append four new bulleted paragraphs after it:
#   - level 0: "RTDE doesn't care ..."
#     - level 1: "Calling "receive" will return ..."
#     - level 1: "Only reason to have more than one key ..."
#     - level 1: "You can try splitting off ..."

# --- New paragraph 1 (list level 0) ---
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.ListFormat.ListOutdent()
$p1.Range.ListFormat.ListOutdent()
$p1.Range.Text = "RTDE doesn’t care about how requested outputs are split up in the XML"

# --- New paragraph 2 (list level 1) ---
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.ListFormat.ListIndent()
$p2.Range.Text = "Calling “receive” will return a packet of data containing all outputs sent via the “send_output_setup” function"

# --- New paragraph 3 (list level 1) ---
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.Text = "Only reason to have more than one key for controller outputs would be for your own organizational purposes"

# --- New paragraph 4 (list level 1) ---
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$p4.Range.Text = "You can try splitting off an output you’re using into a completely different key and the code should work exactly the same "
